# "doc 9 notes started"
# Assignment #6's HW grade (D13) hasn't been graded yet, so its entered
# value is cleared. Also update the Exam 3 score (H8) for student 1.
# All other changed cells (D6, D7, H7, J7) are formulas that recalculate
# automatically from these two edits.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Exam 3 grade for student 1 changes from 0.55 to 0.52
$ws.Range("H8").Value = 0.52

# Clear the recorded HW grade for assignment #6 (row 13) - notes/grading started, no value yet
$ws.Range("D13").ClearContents()

$excel.CalculateFullRebuild()
